$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "time_taken" timestamps on the data sheet (F2:F20) ---
$newTimes = @(
    "2021-10-05 14:21:25.450399",
    "2021-10-05 14:21:25.450408",
    "2021-10-05 14:21:25.450411",
    "2021-10-05 14:21:25.450413",
    "2021-10-05 14:21:25.450416",
    "2021-10-05 14:21:25.450419",
    "2021-10-05 14:21:25.450421",
    "2021-10-05 14:21:25.450424",
    "2021-10-05 14:21:25.450427",
    "2021-10-05 14:21:25.450429",
    "2021-10-05 14:21:25.450432",
    "2021-10-05 14:21:25.450434",
    "2021-10-05 14:21:25.450437",
    "2021-10-05 14:21:25.450439",
    "2021-10-05 14:21:25.450442",
    "2021-10-05 14:21:25.450445",
    "2021-10-05 14:21:25.450448",
    "2021-10-05 14:21:25.450450",
    "2021-10-05 14:21:25.450453"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- Add a new "metadata" worksheet right after "data" ---
$metaSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$metaSheet.Name = "metadata"

# Reuse the same header formatting (bold/border/centered style) used on the
# data sheet's header row by copying it across, then overwrite the values.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$dataSheet.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Copy the A-column style (bold/border/centered, used by A2 on "data") down
# to A2 on "metadata" as well.
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Lipodystrophy - childhood onset"
$metaSheet.Range("C2").Value = 546

# "2.16" must stay a plain text value (not get auto-converted to the number
# 2.16) to match the source data, which stores it as a string. Build it via
# a throwaway text formula and paste only the resulting value back in, which
# keeps the cell as text without attaching a numeric/text style override.
$metaSheet.Range("Z1").Formula = "=""2.""&""16"""
$metaSheet.Range("Z1").Copy()
$metaSheet.Range("D2").PasteSpecial(-4163)
$metaSheet.Range("Z1").ClearContents()

$metaSheet.Range("E2").Value = "2021-07-28T09:58:53.125391Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:21:25.446645"
$metaSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/546/?format=json"

Write-Output "metadata sheet added and timestamps updated"
